$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix gender values (column G): patch bug where "Gender" export always showed "P"
$ws.Range("G2").Value = "Wanita"
$ws.Range("G3").Value = "Pria"
$ws.Range("G4").Value = "Wanita"

# Append a new student row (row 5): Ilham Shiddiq
$ws.Range("A5").Value = "Ilham Shiddiq"

# NISN / NIS / Tanggal Lahir need to stay text (leading zeros / literal date
# string), so force text formatting before assigning, then reset the cell
# style back to the same (unstyled) look as the other data rows.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "0024633245"
$ws.Range("B5").Style = $ws.Range("A4").Style

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "181113836"
$ws.Range("C5").Style = $ws.Range("A4").Style

$ws.Range("D5").Value = "shdqillham123@gmail.com"
$ws.Range("E5").Value = "Cimahi"

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2003-07-03"
$ws.Range("F5").Style = $ws.Range("A4").Style

$ws.Range("G5").Value = "Pria"
$ws.Range("H5").Value = "Islam"
$ws.Range("I5").Value = "Padasuka"
$ws.Range("J5").Value = "SMK1"
$ws.Range("K5").Value = 10
$ws.Range("L5").Value = "Secret"
$ws.Range("M5").Value = "Secret"
$ws.Range("N5").Value = "Secret"
$ws.Range("O5").Value = "12 RPL A"
